$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force Text number format before assignment so Excel
# does not auto-convert numeric-looking strings (e.g. "589.40") into
# floating point numbers, then restore the default "Normal" style so
# no stray style index is left on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.686.37'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.586.07'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.40'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.29%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '187.46'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.574.58'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.623'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.68%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +10.12%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.651'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.70'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000312'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.29%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.57'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.151.36'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.33%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.73'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.679.09'
$ws.Range('D17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.591.45'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.48'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('B20').Style = 'Normal'

$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('C20').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '563.72'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +13.02%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'TRON'
$ws.Range('B21').Style = 'Normal'

$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C21').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.120'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.01'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -5.19%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.70'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +7.43%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '95.81'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.70'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.01%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.30%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.19'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.34'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.38'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.27%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.98%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '65.16'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.67%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.48%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.39'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +8.65%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '564.98'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.422'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +7.26%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '38.19'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0780'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Stacks'
$ws.Range('B42').Style = 'Normal'

$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C42').Style = 'Normal'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.38'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.59%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Maker'
$ws.Range('B43').Style = 'Normal'

$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C43').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.347.01'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.83%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.14%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.17%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0446'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.42'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.20%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.138'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.50'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +23.57%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.998'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.10%  '
$ws.Range('E51').Style = 'Normal'

